$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.023.95"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "1.666.42"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.55"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5097"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2642"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06386"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.88"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07401"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.668.56"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.497"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5826"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008501"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.16"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "26.057.68"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.928"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.74"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.77"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.207"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.11"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.586"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1195"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.61"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06652"
$ws.Range("E28").Value = "  +16.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.314"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.521"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.509"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.633"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.015"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6078"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.204"
$ws.Range("E38").Value = "  +5.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "1.074.16"
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8620"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.51"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "1.815.43"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.23"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.040"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05211"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4287"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.948"
$ws.Range("E51").Value = "  +1.98%  "
